# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" sheet (detail holdings) between the existing
# "总计" (totals) sheet and the "2021-Q3" sheet, and adds a corresponding
# summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: push the existing 2021-Q3 total row down to row 3 and
#    write a brand-new row 2 for 2022-Q3.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Copy row 2 (A:D) down to row 3 so the existing style (A column) moves
# along with the data, then fix up the "index" cell (A3) and restore the
# data that belongs there (2021-Q3 numbers).
$totals.Range("A2:D2").Copy($totals.Range("A3:D3"))
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2021-Q3"
$totals.Range("C3").Value = 2
$totals.Range("D3").Value = 0.02

# Now overwrite row 2 with the new 2022-Q3 totals.
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 4
$totals.Range("D2").Value = 0.1

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" detail sheet. Worksheets.Add() with no
#    arguments drops the new sheet immediately before the active sheet
#    ("2021-Q3"), which lands it exactly between "总计" and "2021-Q3".
# ---------------------------------------------------------------------
$wb.Worksheets.Add() | Out-Null
$new = $wb.Worksheets.Item(2)
$new.Name = "2022-Q3"

# Pull the bold/boxed header style that "总计" already uses (style index
# used by its B1:D1 header cells) and paint it across the full header row
# before writing the header text so every header cell picks it up.
$totals.Range("B1").Copy($new.Range("B1:H1"))
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# Same trick for the row-index column (A2:A5), which also uses that style.
$totals.Range("A2").Copy($new.Range("A2:A5"))
$new.Range("A2").Value = 0
$new.Range("A3").Value = 1
$new.Range("A4").Value = 2
$new.Range("A5").Value = 3

# Columns B, D, E, F and G hold text that looks numeric (fund codes with
# leading zeros, percentages stored as text, ...). Force them to Text
# format first so Excel doesn't silently convert them to numbers / drop
# leading zeros.
$new.Range("B2:B5").NumberFormat = "@"
$new.Range("D2:G5").NumberFormat = "@"

$new.Range("B2").Value = "161127"
$new.Range("C2").Value = "易方达标普生物科技指数（QDII-LOF）人民币"
$new.Range("D2").Value = "3.25"
$new.Range("E2").Value = "94.25"
$new.Range("F2").Value = "1.04"
$new.Range("G2").Value = "0.0338"
$new.Range("H2").Value = 5

$new.Range("B3").Value = "012866"
$new.Range("C3").Value = "易方达标普生物科技指数（QDII-LOF）人民币 C"
$new.Range("D3").Value = "3.25"
$new.Range("E3").Value = "94.25"
$new.Range("F3").Value = "1.04"
$new.Range("G3").Value = "0.0338"
$new.Range("H3").Value = 5

$new.Range("B4").Value = "003720"
$new.Range("C4").Value = "易方达标普生物科技指数（QDII-LOF）美元A"
$new.Range("D4").Value = "3.12"
$new.Range("E4").Value = "94.25"
$new.Range("F4").Value = "1.04"
$new.Range("G4").Value = "0.0324"
$new.Range("H4").Value = 5

$new.Range("B5").Value = "012867"
$new.Range("C5").Value = "易方达标普生物科技指数（QDII-LOF）美元 C"
$new.Range("D5").Value = "0.13"
$new.Range("E5").Value = "94.25"
$new.Range("F5").Value = "1.04"
$new.Range("G5").Value = "0.0014"
$new.Range("H5").Value = 5

# ---------------------------------------------------------------------
# 3) Restore "2021-Q3" as the selected/active sheet (it was the active
#    sheet before this edit, and a brand-new sheet would otherwise steal
#    that state).
# ---------------------------------------------------------------------
$wb.Worksheets.Item(3).Activate()
